$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 87-88 (existing rows 87-101 shift down to 89-103,
# and the former last row 102 becomes row 104). Excel's Insert() inherits
# formatting (including the date number format in column D) from the row
# directly above, matching the target workbook.
$ws.Rows("87:88").Insert()

# New row 87: Poroto verde entry for 2021-11-04, Arica y Parinacota
$ws.Cells.Item(87, 1).Value2 = 5
$ws.Cells.Item(87, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(87, 3).Value2 = "Maule"
$ws.Cells.Item(87, 4).Value2 = 44504
$ws.Cells.Item(87, 5).Value2 = 7
$ws.Cells.Item(87, 6).Value2 = 100112031
$ws.Cells.Item(87, 7).Value2 = "Poroto verde"
$ws.Cells.Item(87, 8).Value2 = "Sin especificar"
$ws.Cells.Item(87, 9).Value2 = "Primera"
$ws.Cells.Item(87, 10).Value2 = 150
$ws.Cells.Item(87, 11).Value2 = 30000
$ws.Cells.Item(87, 12).Value2 = 30000
$ws.Cells.Item(87, 13).Value2 = 30000
$ws.Cells.Item(87, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(87, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(87, 16).Value2 = 1200
$ws.Cells.Item(87, 17).Value2 = 25
$ws.Cells.Item(87, 18).Value2 = "Hortaliza"

# New row 88: Poroto verde entry for 2021-11-04, Región del Maule
$ws.Cells.Item(88, 1).Value2 = 5
$ws.Cells.Item(88, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(88, 3).Value2 = "Maule"
$ws.Cells.Item(88, 4).Value2 = 44504
$ws.Cells.Item(88, 5).Value2 = 7
$ws.Cells.Item(88, 6).Value2 = 100112031
$ws.Cells.Item(88, 7).Value2 = "Poroto verde"
$ws.Cells.Item(88, 8).Value2 = "Sin especificar"
$ws.Cells.Item(88, 9).Value2 = "Primera"
$ws.Cells.Item(88, 10).Value2 = 100
$ws.Cells.Item(88, 11).Value2 = 50000
$ws.Cells.Item(88, 12).Value2 = 50000
$ws.Cells.Item(88, 13).Value2 = 50000
$ws.Cells.Item(88, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(88, 15).Value2 = "Región del Maule"
$ws.Cells.Item(88, 16).Value2 = 2000
$ws.Cells.Item(88, 17).Value2 = 25
$ws.Cells.Item(88, 18).Value2 = "Hortaliza"
